# Apply "measured the geometric info of AMS2" edit.
# Workbook is already open; AMS2 is the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AMS2")

# ---------------------------------------------------------------------------
# 1. Make room for two new header rows at the very top of the sheet
#    (D value + Pixel size) by inserting a single row above row 1.
#    Rows 1.. shift down by one; everything below lines up with the target
#    layout (old row N -> new row N+1).
# ---------------------------------------------------------------------------
$ws.Rows("1:1").Insert()

# The insert leaves row 1 blank/unformatted while row 2 now carries the
# formatting that used to belong to row 1 ("D, mm" cell). Copy that
# formatting back up to row 1 so both header rows share the same style.
$ws.Range("A2:B2").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Fill in the new top rows with the measured geometry reference values.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "D"
$ws.Range("B1").Value = 552

$ws.Range("A2").Value = "Pixel"
$ws.Range("B2").Value = 0.0061

# ---------------------------------------------------------------------------
# 3. Update the column sub-headers ("a, mm"/"h, mm" -> "a"/"2h") for each of
#    the three data blocks.
# ---------------------------------------------------------------------------
function Set-SubHeaders($row) {
    $ws.Range("B$row").Value = "a"
    $ws.Range("C$row").Value = "2h"
    $ws.Range("E$row").Value = "a"
    $ws.Range("F$row").Value = "2h"
    $ws.Range("H$row").Value = "a"
    $ws.Range("I$row").Value = "2h"
    if ($row -ne 33) {
        $ws.Range("K$row").Value = "a"
        $ws.Range("L$row").Value = "2h"
    }
}

Set-SubHeaders 5
Set-SubHeaders 19
Set-SubHeaders 33

# ---------------------------------------------------------------------------
# 4. Fill in the newly measured a / 2h values for each crack (#1..#11) block.
# ---------------------------------------------------------------------------

# --- Block #1 / #2 / #3 / #4  (rows 6-15) ---
$ws.Range("B9").Value  = 24
$ws.Range("C9").Value  = 78
$ws.Range("E9").Value  = 23
$ws.Range("F9").Value  = 102

$ws.Range("B10").Value = 29
$ws.Range("C10").Value = 88
$ws.Range("E10").Value = 26
$ws.Range("F10").Value = 104

$ws.Range("B11").Value = 29
$ws.Range("C11").Value = 90
$ws.Range("E11").Value = 34
$ws.Range("F11").Value = 112

$ws.Range("B12").Value = 30
$ws.Range("C12").Value = 92
$ws.Range("E12").Value = 36
$ws.Range("F12").Value = 116

$ws.Range("B13").Value = 30
$ws.Range("C13").Value = 92
$ws.Range("E13").Value = 37
$ws.Range("F13").Value = 120
$ws.Range("K13").Value = 10
$ws.Range("L13").Value = 35

$ws.Range("B14").Value = 32
$ws.Range("C14").Value = 100
$ws.Range("E14").Value = 37
$ws.Range("F14").Value = 120
$ws.Range("H14").Value = 20
$ws.Range("I14").Value = 42
$ws.Range("K14").Value = 13
$ws.Range("L14").Value = 46

$ws.Range("B15").Value = 32
$ws.Range("C15").Value = 102
$ws.Range("E15").Value = 40
$ws.Range("F15").Value = 128
$ws.Range("H15").Value = 27
$ws.Range("I15").Value = 76
$ws.Range("K15").Value = 17
$ws.Range("L15").Value = 52

# --- Block #5 / #6 / #7 / #8  (rows 20-29) ---
$ws.Range("B22").Value = 10
$ws.Range("C22").Value = 60

$ws.Range("B23").Value = 46
$ws.Range("C23").Value = 280
$ws.Range("E23").Value = 24
$ws.Range("F23").Value = 66

$ws.Range("B24").Value = 65
$ws.Range("C24").Value = 328
$ws.Range("E24").Value = 30
$ws.Range("F24").Value = 104
$ws.Range("K24").Value = 18
$ws.Range("L24").Value = 56

$ws.Range("B25").Value = 70
$ws.Range("C25").Value = 328
$ws.Range("E25").Value = 30
$ws.Range("F25").Value = 106
$ws.Range("K25").Value = 18
$ws.Range("L25").Value = 58

$ws.Range("B26").Value = 72
$ws.Range("C26").Value = 332
$ws.Range("E26").Value = 33
$ws.Range("F26").Value = 106
$ws.Range("K26").Value = 19
$ws.Range("L26").Value = 59

$ws.Range("B27").Value = 88
$ws.Range("C27").Value = 360
$ws.Range("E27").Value = 38
$ws.Range("F27").Value = 110
$ws.Range("H27").Value = 22
$ws.Range("I27").Value = 64
$ws.Range("K27").Value = 26
$ws.Range("L27").Value = 61

$ws.Range("B28").Value = 120
$ws.Range("C28").Value = 400
$ws.Range("E28").Value = 51
$ws.Range("F28").Value = 162
$ws.Range("H28").Value = 29
$ws.Range("I28").Value = 86
$ws.Range("K28").Value = 26
$ws.Range("L28").Value = 62

$ws.Range("B29").Value = 146
$ws.Range("C29").Value = 460
$ws.Range("E29").Value = 64
$ws.Range("F29").Value = 192
$ws.Range("H29").Value = 38
$ws.Range("I29").Value = 120
$ws.Range("K29").Value = 32
$ws.Range("L29").Value = 68

# --- Block #9 / #10 / #11  (rows 34-43) ---
$ws.Range("B37").Value = 21
$ws.Range("C37").Value = 52
$ws.Range("H37").Value = 28
$ws.Range("I37").Value = 66

$ws.Range("B38").Value = 26
$ws.Range("C38").Value = 70
$ws.Range("E38").Value = 28
$ws.Range("F38").Value = 64
$ws.Range("H38").Value = 37
$ws.Range("I38").Value = 79

$ws.Range("B39").Value = 26
$ws.Range("C39").Value = 70
$ws.Range("E39").Value = 28
$ws.Range("F39").Value = 64
$ws.Range("H39").Value = 38
$ws.Range("I39").Value = 99

$ws.Range("B40").Value = 28
$ws.Range("C40").Value = 71
$ws.Range("E40").Value = 28
$ws.Range("F40").Value = 76
$ws.Range("H40").Value = 46
$ws.Range("I40").Value = 106

$ws.Range("B41").Value = 35
$ws.Range("C41").Value = 78
$ws.Range("E41").Value = 32
$ws.Range("F41").Value = 78
$ws.Range("H41").Value = 50
$ws.Range("I41").Value = 110

$ws.Range("B42").Value = 36
$ws.Range("C42").Value = 82
$ws.Range("E42").Value = 37
$ws.Range("F42").Value = 85
$ws.Range("H42").Value = 50
$ws.Range("I42").Value = 118

$ws.Range("B43").Value = 38
$ws.Range("C43").Value = 90
$ws.Range("E43").Value = 39
$ws.Range("F43").Value = 86
$ws.Range("H43").Value = 58
$ws.Range("I43").Value = 118

# ---------------------------------------------------------------------------
# 5. Restore the pane/selection bookkeeping that Excel records for the
#    sheet view.
# ---------------------------------------------------------------------------
$ws.Range("O15").Select()
